$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts dob/date/condition/etc. right by one)
$ws.Range("C1").EntireColumn.Insert()

# Gender values per row, in row order (subject first_name order matches existing rows).
# Using an ordered list of (row, value) pairs rather than a hashtable so the
# underlying shared-string table is built in a fixed, predictable order: all the
# "m"/"f" values are written first (in row order), then the "gender" header last.
$genders = @(
    @(2,  "m"),
    @(3,  "m"),
    @(4,  "f"),
    @(5,  "m"),
    @(6,  "f"),
    @(7,  "m"),
    @(8,  "m"),
    @(9,  "f"),
    @(10, "f"),
    @(11, "m"),
    @(12, "f"),
    @(13, "m"),
    @(14, "m"),
    @(15, "m"),
    @(16, "m"),
    @(17, "m"),
    @(18, "m"),
    @(19, "m"),
    @(20, "f"),
    @(21, "f"),
    @(22, "m"),
    @(23, "m"),
    @(24, "f"),
    @(25, "m"),
    @(26, "m"),
    @(27, "f"),
    @(28, "f")
)

foreach ($pair in $genders) {
    $ws.Cells.Item($pair[0], 3).Value = $pair[1]
}

# Header (added last so "gender" is appended to the shared-string table after "m"/"f")
$ws.Range("C1").Value = "gender"

$ws.Range("F11").Select()
